# Atualização de bases das ligas, do dia: 19-04-2024 às 21:40
#
# This script re-applies a data refresh to the "Venezuela Primera Division"
# sheet:
#   - Five pairs of existing match rows had their underlying records
#     re-matched to a (corrected) id / match, which in practice shows up as
#     a full swap of every data column (B:AC) between the two rows in each
#     pair (column A, the running index, does not move).
#   - The very last match (row 203) got new / corrected data (new id, new
#     kickoff time, final score + full odds panel) and the placeholder row
#     right after it (row 204, a fixture with odds only and no result yet)
#     is removed because that fixture's data now lives in row 203.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param($sheet, $rowA, $rowB)

    $rangeA = $sheet.Range("B$rowA`:AC$rowA")
    $rangeB = $sheet.Range("B$rowB`:AC$rowB")

    $valuesA = $rangeA.Value2
    $valuesB = $rangeB.Value2

    $rangeA.Value = $valuesB
    $rangeB.Value = $valuesA
}

# 1) Row-pair content swaps (column A / the sequence index is left alone)
Swap-RowData $ws 93  98
Swap-RowData $ws 94  95
Swap-RowData $ws 97  99
Swap-RowData $ws 100 101
Swap-RowData $ws 116 117

# 2) Row 203 becomes a finished match with a full data row
$ws.Range("B203").Value = 8111023
$ws.Range("E203").Value = 45400.79166666666
$ws.Range("F203").Value = "Estudiantes Merida"
$ws.Range("G203").Value = "Deportivo Rayo Zuliano"
$ws.Range("H203").Value = 2
$ws.Range("I203").Value = 1
$ws.Range("J203").Value = "H"
$ws.Range("K203").Value = 2.15
$ws.Range("L203").Value = 3.3
$ws.Range("M203").Value = 3
$ws.Range("N203").Value = 1.95
$ws.Range("O203").Value = 3.25
$ws.Range("P203").Value = 3.4
$ws.Range("Q203").Value = -0.5
$ws.Range("R203").Value = 1.975
$ws.Range("S203").Value = 1.825
$ws.Range("T203").Value = 2.5
$ws.Range("U203").Value = 1.925
$ws.Range("V203").Value = 1.875
$ws.Range("W203").Value = 0.95
$ws.Range("X203").Value = -1
$ws.Range("Y203").Value = -1
$ws.Range("Z203").Value = 0.9750000000000001
$ws.Range("AA203").Value = -1
$ws.Range("AB203").Value = 0.925
$ws.Range("AC203").Value = -1

# 3) Old row 204 (now folded into row 203) is removed entirely
$ws.Rows.Item(204).Delete()
